$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02, 1.029984129814187, 1.038358071646851, 0.9926147277508489, 1.04533250382127, 1.037461549807567, 1.035127883980919, 1.041146332803674, 0.9955398523336033, 1.048101040408517, 1.015732391937389),
    @(1.02, 1.030722132024416, 1.038935861409737, 0.9936372048519304, 1.046085321330798, 1.037623489110281, 1.035508229249379, 1.041534733257742, 0.9963617723202692, 1.048665429670427, 1.015857797104481),
    @(1.02, 1.031200124979407, 1.039310090585289, 0.9942998659930995, 1.046573293778808, 1.037727109531837, 1.035754096173783, 1.041785713123771, 0.9968940712668345, 1.049030810517632, 1.015938858267929),
    @(1.02, 1.031401180492373, 1.039467500893316, 0.9945786998346017, 1.046778638495285, 1.03777039173206, 1.035857399406796, 1.041891141953991, 0.997117960005301, 1.049184458690633, 1.01597291573727),
    @(1.02, 1.031434944832062, 1.039493935666745, 0.9946255319796338, 1.046813128505388, 1.037777642569006, 1.035874740963206, 1.041908839003773, 0.9971555583673453, 1.049210259309082, 1.015978632912673),
    @(1.02, 1.031202811073923, 1.03931219358051, 0.9943035907982488, 1.046576036817543, 1.037727688971348, 1.035755476750989, 1.041787122197296, 0.9968970624462087, 1.049032863409027, 1.015939313427179),
    @(1.02, 1.030233445498777, 1.038553262569797, 0.9929600610674301, 1.045586744875026, 1.037516518667396, 1.035256472743432, 1.041277664348729, 0.995817528259106, 1.048291739371719, 1.015774790400787),
    @(1.02, 1.028528878802084, 1.037218777033428, 0.9906006454969559, 1.04385008678885, 1.03713552535637, 1.034375369192627, 1.040377384122891, 0.9939188001724441, 1.046987267441303, 1.015484253385698),
    @(1.02, 1.027395022945905, 1.036331151500093, 0.989033133672735, 1.042696883783981, 1.03687560957245, 1.033786838327378, 1.039775567326298, 0.9926553831429383, 1.046118724588313, 1.015290167486648),
    @(1.02, 1.02690467339374, 1.035947306410004, 0.988355674866747, 1.042198645634249, 1.036761670150622, 1.033531746821613, 1.039514605320246, 0.9921088820399291, 1.045742920870216, 1.015206038190921),
    @(1.02, 1.026722630482689, 1.035804806775892, 0.9881042295826724, 1.042013746279111, 1.036719139405039, 1.033436957686439, 1.039417618064291, 0.9919059725120875, 1.04560337471268, 1.015174775934073),
    @(1.02, 1.026761674986007, 1.035835369890845, 0.9881581567098651, 1.042053400158237, 1.036728271818378, 1.033457291921251, 1.039438424596089, 0.9919494934313052, 1.045633305812993, 1.015181482370338),
    @(1.02, 1.026889623725752, 1.035935525752746, 0.9883348863814464, 1.042183358345927, 1.036758158795581, 1.033523912264979, 1.039506589423739, 0.9920921077337197, 1.045731385031486, 1.01520345430336),
    @(1.02, 1.026968469780773, 1.035997245417395, 0.9884438009545853, 1.042263452266969, 1.036776545537472, 1.033564954422014, 1.039548580865382, 0.9921799884222134, 1.045791820770762, 1.015216990241892),
    @(1.02, 1.027427579009337, 1.036356636795915, 0.9890781214508737, 1.042729973711792, 1.036883142058774, 1.033803762677881, 1.039792878814449, 0.9926916645766087, 1.046143671523073, 1.015295749038335),
    @(1.02, 1.02771573281684, 1.036582209499516, 0.989476357848556, 1.043022908089439, 1.036949634607359, 1.03395349379061, 1.039946022060682, 0.9930127773699352, 1.046364454760343, 1.015345128891346),
    @(1.02, 1.027883867506801, 1.036713830602526, 0.9897087662937556, 1.043193878352392, 1.036988284009076, 1.034040804842288, 1.040035311997428, 0.9932001317071769, 1.046493260882554, 1.015373922723882),
    @(1.02, 1.02794120711474, 1.03675871813897, 0.9897880325774034, 1.043252192805611, 1.037001439602691, 1.034070571445998, 1.040065751420409, 0.9932640239640975, 1.046537184964155, 1.015383739204103),
    @(1.02, 1.027684810477028, 1.036558002666475, 0.9894336180360679, 1.042991467960943, 1.036942514499696, 1.033937431591358, 1.039929594945448, 0.9929783193494215, 1.046340764002873, 1.015339831786738),
    @(1.02, 1.026851943379676, 1.035906030191608, 0.9882828385668249, 1.042145084229839, 1.036749363576167, 1.03350429523917, 1.039486518069624, 0.9920501090198102, 1.045702501914656, 1.015196984468924),
    @(1.02, 1.02632883592198, 1.035496558735204, 0.9875604150241495, 1.041613905060829, 1.036626715635214, 1.033231753014228, 1.039207624327348, 0.9914670000341481, 1.045301457130216, 1.01510709633377),
    @(1.02, 1.026606092372598, 1.035713584023245, 0.9879432794643023, 1.041895399933451, 1.036691847675253, 1.033376252475043, 1.039355500432533, 0.991776070289318, 1.045514033683509, 1.015154754631853),
    @(1.02, 1.02769878276792, 1.036568940543195, 0.9894529299347244, 1.043005674072621, 1.0369457321858, 1.033944689483911, 1.039937017762925, 0.9929938892766442, 1.046351468753249, 1.015342225346932),
    @(1.02, 1.02896911408985, 1.037563423996451, 0.9912096547607049, 1.044298258789675, 1.037235068625136, 1.034603360202214, 1.040610421632954, 0.9944092447426414, 1.047324317803258, 1.015559435656775)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}